$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "877×5=" "256×2="
Replace-Text "622×8=" "538×9="
Replace-Text "985×6=" "867×6="
Replace-Text "269×8=" "316×2="
Replace-Text "971×9=" "215×6="
Replace-Text "220×2=" "572×8="
Replace-Text "143×6=" "607×6="
Replace-Text "456×2=" "737×5="
Replace-Text "293×5=" "564×3="
Replace-Text "175×4=" "877×5="
Replace-Text "563×6=" "443×6="
Replace-Text "481×8=" "911×6="
Replace-Text "849×8=" "847×2="
Replace-Text "728×3=" "469×7="
Replace-Text "244×9=" "953×7="
Replace-Text "515×8=" "356×7="
Replace-Text "139×2=" "622×8="
Replace-Text "108×6=" "189×9="
Replace-Text "118×3=" "619×2="
Replace-Text "899×2=" "909×8="
Replace-Text "331×5=" "745×8="
Replace-Text "541×2=" "143×5="
Replace-Text "561×9=" "784×9="
Replace-Text "447×2=" "565×9="
Replace-Text "586×8=" "838×8="
